$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 2798.6667
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2798.6667
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2798.6667
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9306.6667
# Row 137
$ws.Range("H137").Value = 2949.439
$ws.Range("I137").Value = 2159.1143
$ws.Range("J137").Value = 7559.6665
$ws.Range("K137").Value = 6477.342900000001
$ws.Range("L137").Value = 22678.9995
$ws.Range("M137").Value = -3927.342900000001
$ws.Range("N137").Value = -27778.9995
# Row 138
$ws.Range("H138").Value = 2136.41
$ws.Range("I138").Value = 1662.3928
$ws.Range("K138").Value = 4987.178400000001
$ws.Range("M138").Value = 152.8215999999993
# Row 141
$ws.Range("H141").Value = 4890.357
$ws.Range("I141").Value = 2383.8635
$ws.Range("J141").Value = 14080.833
$ws.Range("K141").Value = 7151.5905
$ws.Range("L141").Value = 42242.499
$ws.Range("M141").Value = -1971.5905
$ws.Range("N141").Value = -52602.499

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2315.6553
$ws.Range("I45").Value = 2152.1667
$ws.Range("K45").Value = 2152.1667
$ws.Range("M45").Value = -1775.1667
# Row 61
$ws.Range("H61").Value = 2470.7778
$ws.Range("I61").Value = 2052.7585
$ws.Range("J61").Value = 3228.4375
$ws.Range("K61").Value = 2052.7585
$ws.Range("L61").Value = 3228.4375
$ws.Range("M61").Value = -1840.7585
$ws.Range("N61").Value = -3652.4375
# Row 74
$ws.Range("H74").Value = 1879.2778
$ws.Range("I74").Value = 1378.5
$ws.Range("K74").Value = 1378.5
$ws.Range("M74").Value = -504.5
# Row 77
$ws.Range("H77").Value = 1879.2778
$ws.Range("I77").Value = 1378.5
$ws.Range("K77").Value = 6892.5
$ws.Range("M77").Value = -2524.5
# Row 132
$ws.Range("H132").Value = 4635.7646
$ws.Range("I132").Value = 3639.3333
$ws.Range("J132").Value = 5422.421
$ws.Range("K132").Value = 10917.9999
$ws.Range("L132").Value = 16267.263
$ws.Range("M132").Value = -8387.999899999999
$ws.Range("N132").Value = -21327.263
# Row 136
$ws.Range("H136").Value = 2470.7778
$ws.Range("I136").Value = 2052.7585
$ws.Range("J136").Value = 3228.4375
$ws.Range("K136").Value = 6158.2755
$ws.Range("L136").Value = 9685.3125
$ws.Range("M136").Value = -3608.2755
$ws.Range("N136").Value = -14785.3125

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3400.5417
$ws.Range("I134").Value = 3710
$ws.Range("J134").Value = 3179.5
$ws.Range("K134").Value = 11130
$ws.Range("L134").Value = 9538.5
$ws.Range("M134").Value = -8595
$ws.Range("N134").Value = -14608.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5347.327
$ws.Range("I31").Value = 1299.2084
$ws.Range("J31").Value = 8481.354499999999
$ws.Range("K31").Value = 1299.2084
$ws.Range("L31").Value = 8481.354499999999
$ws.Range("M31").Value = -1004.2084
$ws.Range("N31").Value = -9071.354499999999
# Row 34
$ws.Range("H34").Value = 5347.327
$ws.Range("I34").Value = 1299.2084
$ws.Range("J34").Value = 8481.354499999999
$ws.Range("K34").Value = 1299.2084
$ws.Range("L34").Value = 8481.354499999999
$ws.Range("M34").Value = -1097.2084
$ws.Range("N34").Value = -8885.354499999999
# Row 58
$ws.Range("H58").Value = 1801.5
$ws.Range("I58").Value = 2504
$ws.Range("J58").Value = 1677.5294
$ws.Range("K58").Value = 2504
$ws.Range("L58").Value = 1677.5294
$ws.Range("M58").Value = -2301
$ws.Range("N58").Value = -2083.5294
# Row 99
$ws.Range("H99").Value = 1742.6364
$ws.Range("I99").Value = 907.6
$ws.Range("J99").Value = 1988.2354
$ws.Range("K99").Value = 907.6
$ws.Range("L99").Value = 1988.2354
$ws.Range("M99").Value = 590.4
$ws.Range("N99").Value = -4984.2354
# Row 122
$ws.Range("H122").Value = 1952.9412
$ws.Range("I122").Value = 1600
$ws.Range("K122").Value = 4800
$ws.Range("M122").Value = -2350
# Row 126
$ws.Range("H126").Value = 1742.6364
$ws.Range("I126").Value = 907.6
$ws.Range("J126").Value = 1988.2354
$ws.Range("K126").Value = 2722.8
$ws.Range("L126").Value = 5964.706200000001
$ws.Range("M126").Value = -252.8000000000002
$ws.Range("N126").Value = -10904.7062
# Row 132
$ws.Range("H132").Value = 9262024
$ws.Range("I132").Value = 2569.4443
$ws.Range("J132").Value = 18521480
$ws.Range("K132").Value = 7708.3329
$ws.Range("L132").Value = 55564440
$ws.Range("M132").Value = -5178.3329
$ws.Range("N132").Value = -55569500
# Row 134
$ws.Range("H134").Value = 2384.5293
$ws.Range("I134").Value = 2524.7778
$ws.Range("J134").Value = 2226.75
$ws.Range("K134").Value = 7574.3334
$ws.Range("L134").Value = 6680.25
$ws.Range("M134").Value = -5039.3334
$ws.Range("N134").Value = -11750.25
# Row 136
$ws.Range("H136").Value = 1801.5
$ws.Range("I136").Value = 2504
$ws.Range("J136").Value = 1677.5294
$ws.Range("K136").Value = 7512
$ws.Range("L136").Value = 5032.5882
$ws.Range("M136").Value = -4962
$ws.Range("N136").Value = -10132.5882

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 7254692.5
$ws.Range("I137").Value = 27797912
$ws.Range("K137").Value = 83393736
$ws.Range("M137").Value = -83388636

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3191
$ws.Range("I132").Value = 3087.4285
$ws.Range("J132").Value = 3256.9092
$ws.Range("K132").Value = 9262.2855
$ws.Range("L132").Value = 9770.7276
$ws.Range("M132").Value = -6732.2855
$ws.Range("N132").Value = -14830.7276

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2937.158
$ws.Range("I132").Value = 1929.4286
$ws.Range("K132").Value = 5788.2858
$ws.Range("M132").Value = -3258.2858
# Row 136
$ws.Range("H136").Value = 13891119
$ws.Range("I136").Value = 2570
$ws.Range("J136").Value = 23811510
$ws.Range("K136").Value = 7710
$ws.Range("L136").Value = 71434530
$ws.Range("M136").Value = -5160
$ws.Range("N136").Value = -71439630

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5053044
$ws.Range("I132").Value = 2943.5715
$ws.Range("J132").Value = 8774171
$ws.Range("K132").Value = 8830.7145
$ws.Range("L132").Value = 26322513
$ws.Range("M132").Value = -6300.7145
$ws.Range("N132").Value = -26327573
# Row 136
$ws.Range("H136").Value = 2744.7942
$ws.Range("I136").Value = 2531.7273
$ws.Range("J136").Value = 3135.4167
$ws.Range("K136").Value = 7595.1819
$ws.Range("L136").Value = 9406.250100000001
$ws.Range("M136").Value = -5045.1819
$ws.Range("N136").Value = -14506.2501
